$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9618077874183655
$ws.Range("B1").Value = 3.171871662139893
$ws.Range("C1").Value = 6.808313846588135
$ws.Range("D1").Value = 1.952882528305054
$ws.Range("E1").Value = 1.362452983856201
